$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24.24049683219419
$ws.Range("C2").Value = 58.37224775547493
$ws.Range("D2").Value = 249.2907120441463
$ws.Range("E2").Value = 84.54063340776146
$ws.Range("F2").Value = 2.58251337832879
